$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns C (Estadístico) and D (p-value) for rows 2-11
# per corrected Diebold-Mariano comparisons (Cap1 revision)

$ws.Range("C2").Value = -0.6850713764914885
$ws.Range("D2").Value = 0.5004539028202741

$ws.Range("C3").Value = -0.6085225578942852
$ws.Range("D3").Value = 0.549072987495677

$ws.Range("C4").Value = 0.275505423592347
$ws.Range("D4").Value = 0.7854991039923827

$ws.Range("C5").Value = -1.578590196599814
$ws.Range("D5").Value = 0.12870182924078

$ws.Range("C6").Value = 0.04517737866197612
$ws.Range("D6").Value = 0.9643734966714357

$ws.Range("C7").Value = 1.039221045708155
$ws.Range("D7").Value = 0.3099893042654869

$ws.Range("C8").Value = -0.9396035856094578
$ws.Range("D8").Value = 0.3576236635296375

$ws.Range("C9").Value = 0.7196923389574351
$ws.Range("D9").Value = 0.4792909589912662

$ws.Range("C10").Value = -0.9190732058929473
$ws.Range("D10").Value = 0.3680249165758775

$ws.Range("C11").Value = -1.572673085304842
$ws.Range("D11").Value = 0.1300667684929449
